$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.109.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.993.02'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.013'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.013'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4973'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4194'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.51'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08945'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.14'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.000.30'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.017'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.432'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.014'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.40'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001106'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06764'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.54'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.978'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.127.89'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.30%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.291'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.232.92'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.83'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.23'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.90%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.27%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '127.29'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.048'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09857'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.528'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.821'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02421'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.320'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.055'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06396'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6493'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.48'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1986'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.012'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.36'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.185'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.496'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.06%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06973'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.38%  '
